$d = $word.ActiveDocument

# --- Config: the text already at the end of the target paragraph, and the
# --- new sentence that needs to be appended as its own run.
$existingText = "晴，今天是高考的第一天，上午考语文，下午考数学。"
$newText = "今天天气不错"

# 1) Locate the run containing $existingText via Find; this repositions the
#    Range's Start/End to the matched text bounds.
$matchRange = $d.Content
$found = $matchRange.Find.Execute($existingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text in document"
}

# 2) Resolve the paragraph that owns the matched run (walk the Paragraphs
#    collection rather than trust Range.Paragraphs, which can come back
#    re-collapsed to the Find hit instead of the whole paragraph).
$targetParagraph = $null
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Start -le $matchRange.Start -and $candidate.Range.End -ge $matchRange.End) {
        $targetParagraph = $candidate
        $targetIndex = $i
        break
    }
}
if ($null -eq $targetParagraph) {
    throw "Could not resolve paragraph for matched text"
}

$paraRange = $targetParagraph.Range

# 3) Pull the paragraph's real OOXML so we keep its exact w:p attributes
#    (w14:paraId/w14:textId/rsid...) and the matched run's rPr verbatim.
$paraOpenXml = $paraRange.WordOpenXML
$bodyTag = "<w:body>"
$bodyIdx = $paraOpenXml.IndexOf($bodyTag)
$afterBody = $paraOpenXml.Substring($bodyIdx + $bodyTag.Length)
$closeTag = "</w:p>"
$closeIdx = $afterBody.IndexOf($closeTag)
$paraXml = $afterBody.Substring(0, $closeIdx + $closeTag.Length)

# Opening <w:p ...> tag, attributes intact.
if ($paraXml -notmatch '^(<w:p\b[^>]*>)') {
    throw "Could not parse paragraph opening tag"
}
$pOpenTag = $Matches[1]

# rPr that decorates the run carrying $existingText (fall back to no rPr).
$runRPr = ""
$rPrMatch = [regex]::Match($paraXml, '<w:r\b[^>]*><w:rPr>(.*?)</w:rPr><w:t\b')
if ($rPrMatch.Success) {
    $runRPr = "<w:rPr>" + $rPrMatch.Groups[1].Value + "</w:rPr>"
}

function Escape-Xml($text) {
    $text = $text -replace '&', '&amp;'
    $text = $text -replace '<', '&lt;'
    $text = $text -replace '>', '&gt;'
    return $text
}

# 4) Build the replacement paragraph: same opening tag, no w:pPr, the
#    original run followed by a brand-new run for $newText carrying the
#    same run formatting.
$run1 = "<w:r>" + $runRPr + "<w:t xml:space=`"preserve`">" + (Escape-Xml $existingText) + "</w:t></w:r>"
$run2 = "<w:r>" + $runRPr + "<w:t xml:space=`"preserve`">" + (Escape-Xml $newText) + "</w:t></w:r>"
$finalXml = $pOpenTag + $run1 + $run2 + "</w:p>"

# 5) Swap the paragraph: drop the old one, insert the freshly-built OOXML
#    right after the previous paragraph (or at the very start of the body
#    if this was the first paragraph).
$prevParagraph = $null
if ($targetIndex -gt 1) {
    $prevParagraph = $d.Paragraphs.Item($targetIndex - 1)
}

$paraRange.Delete()

if ($null -eq $prevParagraph) {
    $insertPoint = $d.Range(0, 0)
} else {
    $insertPoint = $d.Range($prevParagraph.Range.End, $prevParagraph.Range.End)
}
$insertPoint.InsertXML($finalXml)
